# pasar3aB.xlsx — "comprobado para tres numeros"
#
# Hoja1 held a worked example for sorting/comparing 4 numbers (mirroring
# Hoja2). This edit reworks Hoja1 into the 3-number case: the two spacer
# columns collapse (the running-total column merges into the adjacent
# blank column, and one more blank column is dropped before the sorted
# "Resultado en B" / "Objetivos correctos en B" blocks), shifting the
# right-hand blocks left, and the three input columns (D:F) get an
# orange highlight fill. Hoja1 also becomes the active tab/sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")

# --- Collapse the now-unused spacer columns -------------------------------
# Old layout:  ... H(blank) I(blank) J(value+fill) K L M N(labels) ... R S T(sorted+fill) ... Z AA AB AC(objetivos)
# New layout:  ... H(blank) I(value+fill) J K L M(labels)          ... O P Q(sorted+fill) ... W X Y Z (objetivos)
# Deleting column I merges the value into the former second blank column,
# and deleting the following two spacer columns (originally N:O) brings the
# two trailing blocks in by the remaining two columns.
$ws1.Range("I1:I14").Delete(-4159) | Out-Null
$ws1.Range("N1:O14").Delete(-4159) | Out-Null

# --- Highlight the three input-number columns with the new orange fill ---
$ws1.Range("D3:F8").Interior.Color = 49407

# --- View: make Hoja1 the active sheet/tab, update zoom & selection ------
$ws1.Activate()
$excel.ActiveWindow.Zoom = 170
$ws1.Range("H14").Select()

$ws2 = $wb.Worksheets.Item("Hoja2")
$excel.ActiveWindow.Zoom = 170
$ws1.Activate()

# --- Print setup on Hoja1 --------------------------------------------------
$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1
